# Insert a new data row before the current row 91 ("Fecha" 2023-10-04 -> 45084 record
# for the Mango/Terminal Hortofrutícola Agro Chillán series), shifting the existing
# rows 91..206 down to 92..207 and growing the used range to A1:T207.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(91).Insert()

$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = 45225
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100108
$ws.Range("H91").Value = "Tropicales y subtropicales"
$ws.Range("I91").Value = 100108002
$ws.Range("J91").Value = "Mango"
$ws.Range("K91").Value = "Sin especificar"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 100
$ws.Range("N91").Value = 10000
$ws.Range("O91").Value = 10000
$ws.Range("P91").Value = 10000
$ws.Range("Q91").Value = "`$/bandeja 4 kilos"
$ws.Range("R91").Value = "Brasil"
$ws.Range("S91").Value = 2500
$ws.Range("T91").Value = 4
